$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$taskDate = (Get-Date -Year 2021 -Month 5 -Day 4 -Hour 0 -Minute 0 -Second 0).Date

# Set text values in the same order the shared strings were originally authored
# 13: MCD numerique (B10)
$ws.Range("B10").Value = "MCD numerique"
# 14: MLD numerique (B11)
$ws.Range("B11").Value = "MLD numerique"
# 15: MCD-MLD papier (B9)
$ws.Range("B9").Value = "MCD-MLD papier"
# 16: Preparation du template / maquette visuel (B12)
$ws.Range("B12").Value = "Preparation du template / maquette visuel"
# 17: Utilisation d'un logiciel... (E12)
$ws.Range("E12").Value = "Utilisation d'un logiciel nommé nicepage pour crée un template. Il me faut demander a mon chef de projet si il considere ok pour une maquette visuel"
# 18: Arborescence du site (B13)
$ws.Range("B13").Value = "Arborescence du site"
# 19: Documentation (B14)
$ws.Range("B14").Value = "Documentation"
# 20: MCD-MLD decrit + mise en page + premier point de la docs ecrit (E14)
$ws.Range("E14").Value = "MCD-MLD decrit + mise en page + premier point de la docs ecrit"
# 21: Planification initiale (E6)
$ws.Range("E6").Value = "Planification initiale"

# Now fill in the remaining numeric values (dates / durations)
$ws.Range("C9").Value = $taskDate
$ws.Range("D9").Value = 40

$ws.Range("C10").Value = $taskDate
$ws.Range("D10").Value = 25

$ws.Range("C11").Value = $taskDate
$ws.Range("D11").Value = 60

$ws.Range("C12").Value = $taskDate
$ws.Range("D12").Value = 120
$ws.Rows.Item(12).RowHeight = 30

$ws.Range("C13").Value = $taskDate
$ws.Range("D13").Value = 40

$ws.Range("C14").Value = $taskDate
$ws.Range("D14").Value = 120

# Update the selection to B15 (matches last saved cursor position in the file)
$ws.Range("B15").Select()
